# Add a new "res tag" row to the chatbot tag/response table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A41").Value = "เบื่อ"
$ws.Range("B41").Value = "งั้นมาเล่นเกมส์ตอบคำถามกันถ้าตอบถูกหมด 3 ข้อจะได้รางวัลจากน้องบอทแหละ <3 โอเค๊?"

$ws.Range("B41").Select()
